$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column F ("dSF") on specific rows, per repulled data.
$updates = @{
    5  = 2
    7  = 4
    12 = -5
    16 = 6
    17 = 0
    18 = 2
    26 = 9
    28 = 0
    34 = -2
    37 = -3
    38 = 1
    43 = -4
    45 = -2
    46 = 3
    48 = 0
    49 = -1
    51 = -10
    52 = 1
    55 = -3
    58 = -2
    63 = -1
    65 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
